$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column B (rows 2-5) with "NULL" values, matching the new shared
# string entry and the per-row <c r="Bn" t="s"><v>6</v></c> cells added
# in the diff.
$ws.Range("B2").Value = "NULL"
$ws.Range("B3").Value = "NULL"
$ws.Range("B4").Value = "NULL"
$ws.Range("B5").Value = "NULL"

# The diff also moves the active selection from B1 to B5.
$ws.Range("B5").Select()
